$d = $word.ActiveDocument

# The "Summary of Estimated Savings and Implementation Costs" table is the
# first table in the document. Insert a new row before its current first
# row ("Annual Cost Savings") to hold the new "Recommendation Type" / "HVAC"
# pair.
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add($t.Rows.Item(1))
$newRow.Cells.Item(1).Range.Text = "Recommendation Type"
$newRow.Cells.Item(2).Range.Text = "HVAC"
